$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric (e.g. "1.00", "51.815.84").
# Force the whole Price column to Text format first so Excel keeps the exact
# literal strings (matching the original inlineStr cell values) instead of
# auto-converting them to numbers, then restore the default "Normal" style so
# no extra style index is left attached to the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.815.84"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.818.97"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "355.46"
$ws.Range("E5").Value = "  +3.53%  "
$ws.Range("D6").Value = "111.67"
$ws.Range("E6").Value = "  -3.18%  "
$ws.Range("D7").Value = "0.566"
$ws.Range("E7").Value = "  +3.44%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.597"
$ws.Range("E9").Value = "  +3.54%  "
$ws.Range("D10").Value = "40.73"
$ws.Range("E10").Value = "  -4.39%  "
$ws.Range("D11").Value = "0.0856"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "19.89"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "7.78"
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("D15").Value = "3.258.16"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "2.835.99"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("D17").Value = "0.919"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").Value = "51.762.74"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "7.56"
$ws.Range("E19").Value = "  +7.72%  "
$ws.Range("D20").Value = "3.13"
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("D21").Value = "13.37"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").Value = "0.0₃0993"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").Value = "70.09"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "267.60"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").Value = "2.80"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").Value = "26.96"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "10.29"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("B29").Value = "VeChain"
$ws.Range("C29").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D29").Value = "0.0496"
$ws.Range("E29").Value = "  +26.99%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.21"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").Value = "  +5.06%  "
$ws.Range("D33").Value = "34.73"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("E34").Value = "  +2.54%  "
$ws.Range("D35").Value = "5.46"
$ws.Range("E35").Value = "  +10.25%  "
$ws.Range("E36").Value = "  +3.05%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("D40").Value = "18.36"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("E42").Value = "  -4.23%  "
$ws.Range("D43").Value = "23.17"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").Value = "125.89"
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D46").Value = "2.097.29"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").Value = "3.34"
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").Value = "6.01"
$ws.Range("E49").Value = "  +8.34%  "
$ws.Range("D50").Value = "0.986"
$ws.Range("E50").Value = "  +9.72%  "
$ws.Range("D51").Value = "9.05"
$ws.Range("E51").Value = "  +2.07%  "

$ws.Range("D2:D51").Style = "Normal"
